$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values are stored as text, matching the source data
# (values like "182.71" would otherwise be auto-converted to numbers by Excel)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.612.03"
$ws.Range("E2").Value = "  -6.99%  "
$ws.Range("D3").Value = "3.297.48"
$ws.Range("E3").Value = "  -7.78%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "182.71"
$ws.Range("E5").Value = "  -11.71%  "
$ws.Range("D6").Value = "522.49"
$ws.Range("E6").Value = "  -6.97%  "
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "3.295.86"
$ws.Range("E8").Value = "  -7.69%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "0.620"
$ws.Range("E10").Value = "  -7.81%  "
$ws.Range("D11").Value = "58.13"
$ws.Range("E11").Value = "  -8.06%  "
$ws.Range("E12").Value = "  -9.75%  "
$ws.Range("E13").Value = "  -8.26%  "
$ws.Range("D14").Value = "9.13"
$ws.Range("E14").Value = "  -9.55%  "
$ws.Range("D15").Value = "3.816.09"
$ws.Range("E15").Value = "  -7.99%  "
$ws.Range("E16").Value = "  -5.41%  "
$ws.Range("D17").Value = "3.295.56"
$ws.Range("D18").Value = "17.76"
$ws.Range("E18").Value = "  -6.90%  "
$ws.Range("D19").Value = "63.552.24"
$ws.Range("E19").Value = "  -6.78%  "
$ws.Range("D20").Value = "10.98"
$ws.Range("E20").Value = "  -9.48%  "
$ws.Range("D21").Value = "0.951"
$ws.Range("E21").Value = "  -10.21%  "
$ws.Range("D22").Value = "371.34"
$ws.Range("E22").Value = "  -6.58%  "
$ws.Range("D23").Value = "11.24"
$ws.Range("E23").Value = "  -9.21%  "
$ws.Range("E24").Value = "  -10.49%  "
$ws.Range("D25").Value = "80.18"
$ws.Range("E25").Value = "  -4.64%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "5.99"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").Value = "2.64"
$ws.Range("E28").Value = "  -7.74%  "
$ws.Range("D29").Value = "11.40"
$ws.Range("E29").Value = "  -7.92%  "
$ws.Range("D30").Value = "8.32"
$ws.Range("E30").Value = "  -8.40%  "
$ws.Range("D31").Value = "651.80"
$ws.Range("E31").Value = "  -8.78%  "
$ws.Range("D32").Value = "28.59"
$ws.Range("E32").Value = "  -8.72%  "
$ws.Range("D33").Value = "6.69"
$ws.Range("E33").Value = "  -11.03%  "
$ws.Range("D34").Value = "11.18"
$ws.Range("E34").Value = "  -7.04%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.106"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "59.25"
$ws.Range("E36").Value = "  -6.91%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -6.98%  "
$ws.Range("D39").Value = "36.16"
$ws.Range("E39").Value = "  -12.16%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "2.997.80"
$ws.Range("E41").Value = "  -4.84%  "
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").Value = "0.0₃0652"
$ws.Range("E43").Value = "  -10.63%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "2.69"
$ws.Range("E44").Value = "  -16.23%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("D46").Value = "0.0389"
$ws.Range("E46").Value = "  -5.05%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.82"
$ws.Range("E48").Value = "  +4.53%  "
$ws.Range("E49").Value = "  -3.85%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.50"
$ws.Range("E50").Value = "  -20.01%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "2.91"
$ws.Range("E51").Value = "  -5.07%  "
